$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb1"
$ws.Cells.Item(2, 3).Value = "Ephb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 10.31211433333333
$ws.Cells.Item(2, 8).Value = 30.936343
$ws.Cells.Item(2, 9).Value = 0.633340936097251
$ws.Cells.Item(2, 10).Value = 0.633340936097251
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.2721246666666666
$ws.Cells.Item(2, 14).Value = 0.8163739999999999
$ws.Cells.Item(2, 15).Value = 0.1154759865526449
$ws.Cells.Item(2, 16).Value = 0.1154759865526449
$ws.Cells.Item(2, 17).Value = 2.806180675586889
$ws.Cells.Item(2, 18).Value = 25.255626080282
$ws.Cells.Item(2, 19).Value = 0.07313566942000566
$ws.Cells.Item(2, 20).Value = 0.07313566942000566

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb1"
$ws.Cells.Item(3, 3).Value = "Ephb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10.31211433333333
$ws.Cells.Item(3, 8).Value = 30.936343
$ws.Cells.Item(3, 9).Value = 0.633340936097251
$ws.Cells.Item(3, 10).Value = 0.633340936097251
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.8319233333333332
$ws.Cells.Item(3, 14).Value = 2.49577
$ws.Cells.Item(3, 15).Value = 0.3530263126440755
$ws.Cells.Item(3, 16).Value = 0.3530263126440755
$ws.Cells.Item(3, 17).Value = 8.578888529901111
$ws.Cells.Item(3, 18).Value = 77.20999676910999
$ws.Cells.Item(3, 19).Value = 0.2235860153169595
$ws.Cells.Item(3, 20).Value = 0.2235860153169595

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb1"
$ws.Cells.Item(4, 3).Value = "Ephb6"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 10.31211433333333
$ws.Cells.Item(4, 8).Value = 30.936343
$ws.Cells.Item(4, 9).Value = 0.633340936097251
$ws.Cells.Item(4, 10).Value = 0.633340936097251
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.252499666666667
$ws.Cells.Item(4, 14).Value = 3.757499
$ws.Cells.Item(4, 15).Value = 0.5314977008032796
$ws.Cells.Item(4, 16).Value = 0.5314977008032796
$ws.Cells.Item(4, 17).Value = 12.91591976512856
$ws.Cells.Item(4, 18).Value = 116.243277886157
$ws.Cells.Item(4, 19).Value = 0.3366192513602858
$ws.Cells.Item(4, 20).Value = 0.3366192513602858

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb1"
$ws.Cells.Item(5, 3).Value = "Ephb6"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.103438
$ws.Cells.Item(5, 8).Value = 12.310314
$ws.Cells.Item(5, 9).Value = 0.2520215719230645
$ws.Cells.Item(5, 10).Value = 0.2520215719230645
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2721246666666666
$ws.Cells.Item(5, 14).Value = 0.8163739999999999
$ws.Cells.Item(5, 15).Value = 0.1154759865526449
$ws.Cells.Item(5, 16).Value = 0.1154759865526449
$ws.Cells.Item(5, 17).Value = 1.116646697937333
$ws.Cells.Item(5, 18).Value = 10.049820281436
$ws.Cells.Item(5, 19).Value = 0.02910243965036422
$ws.Cells.Item(5, 20).Value = 0.02910243965036421

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb1"
$ws.Cells.Item(6, 3).Value = "Ephb6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.103438
$ws.Cells.Item(6, 8).Value = 12.310314
$ws.Cells.Item(6, 9).Value = 0.2520215719230645
$ws.Cells.Item(6, 10).Value = 0.2520215719230645
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8319233333333332
$ws.Cells.Item(6, 14).Value = 2.49577
$ws.Cells.Item(6, 15).Value = 0.3530263126440755
$ws.Cells.Item(6, 16).Value = 0.3530263126440755
$ws.Cells.Item(6, 17).Value = 3.413745819086666
$ws.Cells.Item(6, 18).Value = 30.72371237177999
$ws.Cells.Item(6, 19).Value = 0.08897024624276312
$ws.Cells.Item(6, 20).Value = 0.0889702462427631

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb1"
$ws.Cells.Item(7, 3).Value = "Ephb6"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.103438
$ws.Cells.Item(7, 8).Value = 12.310314
$ws.Cells.Item(7, 9).Value = 0.2520215719230645
$ws.Cells.Item(7, 10).Value = 0.2520215719230645
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.252499666666667
$ws.Cells.Item(7, 14).Value = 3.757499
$ws.Cells.Item(7, 15).Value = 0.5314977008032796
$ws.Cells.Item(7, 16).Value = 0.5314977008032796
$ws.Cells.Item(7, 17).Value = 5.139554727187334
$ws.Cells.Item(7, 18).Value = 46.25599254468599
$ws.Cells.Item(7, 19).Value = 0.1339488860299372
$ws.Cells.Item(7, 20).Value = 0.1339488860299371

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efnb1"
$ws.Cells.Item(8, 3).Value = "Ephb6"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.866538
$ws.Cells.Item(8, 8).Value = 5.599614
$ws.Cells.Item(8, 9).Value = 0.1146374919796846
$ws.Cells.Item(8, 10).Value = 0.1146374919796846
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.2721246666666666
$ws.Cells.Item(8, 14).Value = 0.8163739999999999
$ws.Cells.Item(8, 15).Value = 0.1154759865526449
$ws.Cells.Item(8, 16).Value = 0.1154759865526449
$ws.Cells.Item(8, 17).Value = 0.5079310310706666
$ws.Cells.Item(8, 18).Value = 4.571379279635999
$ws.Cells.Item(8, 19).Value = 0.01323787748227499
$ws.Cells.Item(8, 20).Value = 0.01323787748227499

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efnb1"
$ws.Cells.Item(9, 3).Value = "Ephb6"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.866538
$ws.Cells.Item(9, 8).Value = 5.599614
$ws.Cells.Item(9, 9).Value = 0.1146374919796846
$ws.Cells.Item(9, 10).Value = 0.1146374919796846
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.8319233333333332
$ws.Cells.Item(9, 14).Value = 2.49577
$ws.Cells.Item(9, 15).Value = 0.3530263126440755
$ws.Cells.Item(9, 16).Value = 0.3530263126440755
$ws.Cells.Item(9, 17).Value = 1.552816514753333
$ws.Cells.Item(9, 18).Value = 13.97534863278
$ws.Cells.Item(9, 19).Value = 0.04047005108435282
$ws.Cells.Item(9, 20).Value = 0.04047005108435282

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efnb1"
$ws.Cells.Item(10, 3).Value = "Ephb6"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.866538
$ws.Cells.Item(10, 8).Value = 5.599614
$ws.Cells.Item(10, 9).Value = 0.1146374919796846
$ws.Cells.Item(10, 10).Value = 0.1146374919796846
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.252499666666667
$ws.Cells.Item(10, 14).Value = 3.757499
$ws.Cells.Item(10, 15).Value = 0.5314977008032796
$ws.Cells.Item(10, 16).Value = 0.5314977008032796
$ws.Cells.Item(10, 17).Value = 2.337838222820667
$ws.Cells.Item(10, 18).Value = 21.040544005386
$ws.Cells.Item(10, 19).Value = 0.06092956341305677
$ws.Cells.Item(10, 20).Value = 0.06092956341305676
